$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 29 (shifts existing rows 29-52 down to 32-55)
$ws.Rows("29:31").Insert()

# --- New row 29: Especial ---
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44629
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107011
$ws.Range("J29").Value = "Tuna"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 500
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 13500
$ws.Range("Q29").Value = "$/caja 18 kilos"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 750
$ws.Range("T29").Value = 18

# --- New row 30: Primera ---
$ws.Range("A30").Value = 2
$ws.Range("B30").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44629
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107011
$ws.Range("J30").Value = "Tuna"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 400
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 11000
$ws.Range("P30").Value = 10500
$ws.Range("Q30").Value = "$/caja 18 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 583
$ws.Range("T30").Value = 18

# --- New row 31: Segunda ---
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44629
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100107
$ws.Range("H31").Value = "Otros"
$ws.Range("I31").Value = 100107011
$ws.Range("J31").Value = "Tuna"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 400
$ws.Range("N31").Value = 7000
$ws.Range("O31").Value = 8000
$ws.Range("P31").Value = 7500
$ws.Range("Q31").Value = "$/caja 18 kilos"
$ws.Range("R31").Value = "Provincia de Limarí"
$ws.Range("S31").Value = 417
$ws.Range("T31").Value = 18
